$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text content would otherwise be auto-parsed as a number
# by Excel (losing the literal formatting, e.g. trailing zeros or the
# thousands-dot style used by this sheet) are forced to Text format first.

$ws.Range("D2").Value = "30.307.52"
$ws.Range("E2").Value = "  -0.18%  "

$ws.Range("D3").Value = "1.858.62"
$ws.Range("E3").Value = "  -0.82%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  +0.29%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "232.47"
$ws.Range("E5").Value = "  -2.46%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.003"
$ws.Range("E6").Value = "  +0.26%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4748"
$ws.Range("E7").Value = "  -0.70%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2752"
$ws.Range("E8").Value = "  -2.59%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06439"
$ws.Range("E9").Value = "  -1.33%  "

$ws.Range("D10").Value = "1.852.83"
$ws.Range("E10").Value = "  -1.14%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07411"
$ws.Range("E11").Value = "  -0.62%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "16.02"
$ws.Range("E12").Value = "  -3.14%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.994"
$ws.Range("E13").Value = "  -2.06%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "85.83"
$ws.Range("E14").Value = "  -2.61%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6307"
$ws.Range("E15").Value = "  -3.67%  "

$ws.Range("D16").Value = "30.315.92"
$ws.Range("E16").Value = "  -0.07%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.001"
$ws.Range("E17").Value = "  +0.16%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.83"
$ws.Range("E18").Value = "  -3.67%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "227.09"
$ws.Range("E19").Value = "  +3.84%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007311"
$ws.Range("E20").Value = "  -3.86%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.005"
$ws.Range("E21").Value = "  +0.46%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.087"
$ws.Range("E22").Value = "  -4.11%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.051"
$ws.Range("E23").Value = "  -2.47%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "167.25"
$ws.Range("E24").Value = "  +0.02%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.218"
$ws.Range("E25").Value = "  -1.29%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "17.83"
$ws.Range("E26").Value = "  -3.41%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.870"
$ws.Range("E27").Value = "  -5.11%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.1035"
$ws.Range("E28").Value = "  +10.77%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.383"
$ws.Range("E29").Value = "  -5.17%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.206"
$ws.Range("E30").Value = "  -2.56%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.894"
$ws.Range("E31").Value = "  -3.48%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.04911"
$ws.Range("E32").Value = "  -3.27%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.152"
$ws.Range("E33").Value = "  -4.35%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7244"
$ws.Range("E34").Value = "  -3.47%  "

$ws.Range("E35").Value = "  +0.09%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.717"
$ws.Range("E36").Value = "  +0.07%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01892"
$ws.Range("E37").Value = "  +3.56%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.628"
$ws.Range("E38").Value = "  +0.50%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.9055"
$ws.Range("E39").Value = "  -0.02%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.980"
$ws.Range("E40").Value = "  -4.10%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "105.17"
$ws.Range("E41").Value = "  -1.64%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9961"
$ws.Range("E42").Value = "  -0.76%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4095"
$ws.Range("E43").Value = "  -4.17%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.535"
$ws.Range("E44").Value = "  -6.11%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.024"
$ws.Range("E45").Value = "  -4.67%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "60.88"
$ws.Range("E46").Value = "  -5.34%  "

$ws.Range("B47").Value = "Algorand"
$ws.Range("C47").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1201"
$ws.Range("E47").Value = "  -6.39%  "

$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.763"
$ws.Range("E48").Value = "  -1.85%  "

$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05620"
$ws.Range("E49").Value = "  -0.10%  "

$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.396"
$ws.Range("E50").Value = "  -5.06%  "

$ws.Range("B51").Value = "Elrond"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "32.92"
$ws.Range("E51").Value = "  -2.03%  "
